$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.206.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.62%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.085.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +8.93%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "

# Row 6
$ws.Range("E6").Value = "  -4.15%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.88"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.28%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.90"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.26%  "

# Row 10
$ws.Range("E10").Value = "  +0.05%  "

# Row 12
$ws.Range("E12").Value = "  +6.13%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.72%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.390.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.02%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.838"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.80%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.082.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.79%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.09%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.142.98"
$ws.Range("D18").Style = "Normal"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.33%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0828"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.36%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.26%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "241.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.13%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.65%  "

# Row 24
$ws.Range("E24").Value = "  +0.18%  "

# Row 25
$ws.Range("E25").Value = "  +2.14%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.54%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.20%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.36%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.98%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +54.30%  "

# Row 31
$ws.Range("E31").Value = "  -4.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +21.93%  "

# Row 33
$ws.Range("E33").Value = "  -1.05%  "

# Row 34
$ws.Range("E34").Value = "  +0.13%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0937"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.87%  "

# Row 36
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.11%  "

# Row 37
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +16.84%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.12%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.34%  "

# Row 40
$ws.Range("E40").Value = "  -9.39%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.51%  "

# Row 42
$ws.Range("E42").Value = "  -0.99%  "

# Row 43
$ws.Range("E43").Value = "  +5.79%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.23%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.64%  "

# Row 46
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0870"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.39%  "

# Row 47
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.335.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.15%  "

# Row 48
$ws.Range("E48").Value = "  +4.85%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.79%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.288.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.52%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.31%  "

Write-Output "done"